$d = $word.ActiveDocument

# The document starts with two paragraphs that each contain a single
# inline picture (InlineShape). The commit swaps out the portfolio
# screenshot: both existing pictures are removed, leaving the two
# paragraph marks behind as empty paragraphs. The following paragraph
# (with the manual line breaks) is left untouched.

while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}
